$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Prophet MAE")

$ws.Range("C2").Value = 0.4147712986052281
$ws.Range("C3").Value = 0.09999565895742606
$ws.Range("C4").Value = 0.2021875542765255
$ws.Range("C6").Value = 0.1029787273645813
$ws.Range("C7").Value = 0.2350353010065291
$ws.Range("C8").Value = 0.1326798583002398
$ws.Range("C9").Value = 0.3293286587327779
$ws.Range("C10").Value = 0.1755939442914869
